# CryCompanywiseStockReport_1.xlsx - stock-count correction pass.
#
# The source report lists, per supplying company, a block of stock items
# (A=serial#, B=item code, C=description, D=cost rate, E=sale rate,
#  F=quantity on hand, G=stock value = D*F), each block closed by a
# "Sub Total:" row whose B column is SUM(G) for that block, and the whole
# sheet closed by two rows: a "Sub Total:" (row 1051) that totals every
# per-company subtotal, and a "Grand Total:" (row 1052) mirroring it.
#
# This pass revises the on-hand quantity (F) for a batch of items (mostly
# -1/-2 count corrections, a couple of bigger re-counts, and one item
# zeroed out), restores two item rows that had gotten their code/name/
# rate/qty swapped between each other, and fixes one pair of rows that had
# their item codes swapped. All downstream value/subtotal/grand-total
# cells are then recomputed to stay consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-Num($addr) {
    $v = $ws.Range($addr).Value2
    if ($v -is [double]) { return $v }
    if ($v -is [int]) { return [double]$v }
    if ($v -is [int32]) { return [double]$v }
    if ($v -is [int64]) { return [double]$v }
    # Blank cells (NULL) and any stray text (e.g. the GST-number header
    # row) are not part of the quantity/value columns' running total.
    return 0.0
}

function Set-Num($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Get-Txt($addr) {
    $v = $ws.Range($addr).Value2
    if ($v -eq $null) { return "" }
    return $v
}

function Set-Txt($addr, $val) {
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------
# 1) Two item rows (199/200) had swapped item codes (column B) - restore.
# ---------------------------------------------------------------------
$b199 = Get-Num "B199"
$b200 = Get-Num "B200"
Set-Num "B199" $b200
Set-Num "B200" $b199

# ---------------------------------------------------------------------
# 2) Two item rows (674/675) had their entire records (code, description,
#    cost rate, sale rate, qty, value) swapped with each other - restore
#    by swapping B/C/D/E/F/G between the two rows (serial # in A is left
#    untouched, it belongs to the row position, not the item).
# ---------------------------------------------------------------------
foreach ($col in @("B", "C", "D", "E", "F", "G")) {
    $addr1 = "$col" + "674"
    $addr2 = "$col" + "675"
    $v1 = Get-Txt $addr1
    $v2 = Get-Txt $addr2
    Set-Txt $addr1 $v2
    Set-Txt $addr2 $v1
}

# ---------------------------------------------------------------------
# 3) Quantity-on-hand (F) corrections: row -> new quantity. Value (G) is
#    recomputed as cost-rate (D) * quantity (F) right after, matching how
#    every other row in the sheet is derived.
# ---------------------------------------------------------------------
$qtyFixes = @{
    26 = 55
    32 = 14
    33 = 17
    45 = 33
    52 = 31
    55 = 175
    61 = 40
    98 = 0
    123 = 173
    148 = 86
    153 = 25
    164 = 74
    174 = 5
    208 = 21
    224 = 13
    236 = 22
    266 = 51
    271 = 29
    344 = 31
    412 = 96
    420 = 24
    440 = 43
    450 = 187
    451 = 55
    516 = 21
    522 = 69
    552 = 165
    557 = 128
    597 = 684
    600 = 134
    648 = 80
    658 = 150
    672 = 127
    673 = 30
    674 = 1
    675 = 0
    680 = 55
    681 = 47
    685 = 24
    727 = 26
    734 = 108
    745 = 28
    813 = 39
    906 = 87
    908 = 41
    939 = 50
    947 = 188
    957 = 23
    961 = 42
    965 = 16
    974 = 458
    975 = 334
    977 = 336
    981 = 121
    999 = 9
    1003 = 10
}

foreach ($row in $qtyFixes.Keys) {
    $newQty = $qtyFixes[$row]
    $rate = Get-Num "D$row"
    Set-Num "F$row" $newQty
    Set-Num "G$row" ($rate * $newQty)
}

# ---------------------------------------------------------------------
# 4) Recompute every "Sub Total:" block (SUM of G since the previous
#    block boundary) and the trailing Sub Total/Grand Total pair that
#    totals all the per-company subtotals. Walk the whole used range once.
# ---------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

$runningSum = 0.0
$prevLabel = ""
$subtotalSum = 0.0

for ($r = 1; $r -le $lastRow; $r++) {
    $label = Get-Txt "A$r"

    if ($label -eq "Sub Total:") {
        if ($prevLabel -eq "Sub Total:") {
            # The grand-aggregating "Sub Total:" row: total of every
            # per-company subtotal seen so far.
            Set-Num "B$r" $subtotalSum
        } else {
            Set-Num "B$r" $runningSum
            $subtotalSum = $subtotalSum + $runningSum
        }
        $runningSum = 0.0
    } elseif ($label -eq "Grand Total:") {
        Set-Num "B$r" $subtotalSum
    } else {
        $g = Get-Num "G$r"
        $runningSum = $runningSum + $g
    }

    $prevLabel = $label
}
